$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 511; existing rows 511..612 shift down to 512..613
$ws.Rows.Item(511).Insert()

# Populate the new row 511 with data
$ws.Cells.Item(511, 1).Value2 = 8
$ws.Cells.Item(511, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(511, 3).Value2 = "Coquimbo"
$ws.Cells.Item(511, 4).Value2 = 45209
$ws.Cells.Item(511, 5).Value2 = 4
$ws.Cells.Item(511, 6).Value2 = 100114013
$ws.Cells.Item(511, 7).Value2 = "Zanahoria"
$ws.Cells.Item(511, 8).Value2 = "Sin especificar"
$ws.Cells.Item(511, 9).Value2 = "Primera"
$ws.Cells.Item(511, 10).Value2 = 600
$ws.Cells.Item(511, 11).Value2 = 5800
$ws.Cells.Item(511, 12).Value2 = 6000
$ws.Cells.Item(511, 13).Value2 = 5900
$ws.Cells.Item(511, 14).Value2 = "`$/saco 20 kilos"
$ws.Cells.Item(511, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(511, 16).Value2 = 295
$ws.Cells.Item(511, 17).Value2 = 20
$ws.Cells.Item(511, 18).Value2 = "Hortaliza"

# Apply the same date number format as the other date cells in column D
$ws.Cells.Item(511, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
